# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos sheet
# with refreshed values, as published by the "Updated cryptos list" GitHub
# Action.
#
# The cell values are numeric-looking text (e.g. "211.30", "27.465.80",
# "  -1.70%  ") that must stay stored as literal text, exactly as in the
# original file, instead of being auto-converted to numbers by Excel. To
# achieve this we assign each value with a leading apostrophe (Excel's
# "treat as text" quote-prefix convention) and then reset the cell Style
# back to "Normal" so no stray number formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.Value = "'" + '27.465.80'
$r.Style = "Normal"
$r = $ws.Range('D3')
$r.Value = "'" + '1.617.13'
$r.Style = "Normal"
$r = $ws.Range('E3')
$r.Value = "'" + '  -1.70%  '
$r.Style = "Normal"
$r = $ws.Range('E4')
$r.Value = "'" + '  +0.02%  '
$r.Style = "Normal"
$r = $ws.Range('D5')
$r.Value = "'" + '211.30'
$r.Style = "Normal"
$r = $ws.Range('E5')
$r.Value = "'" + '  -1.02%  '
$r.Style = "Normal"
$r = $ws.Range('E6')
$r.Value = "'" + '  -1.31%  '
$r.Style = "Normal"
$r = $ws.Range('E7')
$r.Value = "'" + '  +0.01%  '
$r.Style = "Normal"
$r = $ws.Range('D8')
$r.Value = "'" + '22.85'
$r.Style = "Normal"
$r = $ws.Range('E9')
$r.Value = "'" + '  +0.21%  '
$r.Style = "Normal"
$r = $ws.Range('E10')
$r.Value = "'" + '  -0.36%  '
$r.Style = "Normal"
$r = $ws.Range('D11')
$r.Value = "'" + '0.0886'
$r.Style = "Normal"
$r = $ws.Range('E11')
$r.Value = "'" + '  -0.62%  '
$r.Style = "Normal"
$r = $ws.Range('D12')
$r.Value = "'" + '1.844.41'
$r.Style = "Normal"
$r = $ws.Range('E12')
$r.Value = "'" + '  -1.78%  '
$r.Style = "Normal"
$r = $ws.Range('D13')
$r.Value = "'" + '1.630.56'
$r.Style = "Normal"
$r = $ws.Range('E13')
$r.Value = "'" + '  -0.91%  '
$r.Style = "Normal"
$r = $ws.Range('E14')
$r.Value = "'" + '  -0.20%  '
$r.Style = "Normal"
$r = $ws.Range('E15')
$r.Value = "'" + '  -2.71%  '
$r.Style = "Normal"
$r = $ws.Range('D16')
$r.Value = "'" + '64.87'
$r.Style = "Normal"
$r = $ws.Range('E16')
$r.Value = "'" + '  +0.93%  '
$r.Style = "Normal"
$r = $ws.Range('D17')
$r.Value = "'" + '27.451.49'
$r.Style = "Normal"
$r = $ws.Range('E17')
$r.Value = "'" + '  -0.83%  '
$r.Style = "Normal"
$r = $ws.Range('D18')
$r.Value = "'" + '231.04'
$r.Style = "Normal"
$r = $ws.Range('E18')
$r.Value = "'" + '  -0.33%  '
$r.Style = "Normal"
$r = $ws.Range('E19')
$r.Value = "'" + '  -1.07%  '
$r.Style = "Normal"
$r = $ws.Range('D20')
$r.Value = "'" + '7.54'
$r.Style = "Normal"
$r = $ws.Range('E20')
$r.Value = "'" + '  -1.94%  '
$r.Style = "Normal"
$r = $ws.Range('E21')
$r.Value = "'" + '  +0.10%  '
$r.Style = "Normal"
$r = $ws.Range('E22')
$r.Value = "'" + '  -0.63%  '
$r.Style = "Normal"
$r = $ws.Range('D23')
$r.Value = "'" + '10.15'
$r.Style = "Normal"
$r = $ws.Range('E23')
$r.Value = "'" + '  +0.43%  '
$r.Style = "Normal"
$r = $ws.Range('E24')
$r.Value = "'" + '  +6.00%  '
$r.Style = "Normal"
$r = $ws.Range('D25')
$r.Value = "'" + '150.87'
$r.Style = "Normal"
$r = $ws.Range('E25')
$r.Value = "'" + '  +0.56%  '
$r.Style = "Normal"
$r = $ws.Range('E26')
$r.Value = "'" + '  -1.85%  '
$r.Style = "Normal"
$r = $ws.Range('D27')
$r.Value = "'" + '0.112'
$r.Style = "Normal"
$r = $ws.Range('E27')
$r.Value = "'" + '  -0.94%  '
$r.Style = "Normal"
$r = $ws.Range('E28')
$r.Value = "'" + '  +0.01%  '
$r.Style = "Normal"
$r = $ws.Range('D29')
$r.Value = "'" + '15.54'
$r.Style = "Normal"
$r = $ws.Range('E29')
$r.Value = "'" + '  -0.89%  '
$r.Style = "Normal"
$r = $ws.Range('E30')
$r.Value = "'" + '  -0.95%  '
$r.Style = "Normal"
$r = $ws.Range('E31')
$r.Value = "'" + '  -1.01%  '
$r.Style = "Normal"
$r = $ws.Range('E32')
$r.Value = "'" + '  -1.28%  '
$r.Style = "Normal"
$r = $ws.Range('D33')
$r.Value = "'" + '1.467.68'
$r.Style = "Normal"
$r = $ws.Range('E33')
$r.Value = "'" + '  +1.57%  '
$r.Style = "Normal"
$r = $ws.Range('E34')
$r.Value = "'" + '  -2.90%  '
$r.Style = "Normal"
$r = $ws.Range('D35')
$r.Value = "'" + '1.54'
$r.Style = "Normal"
$r = $ws.Range('E35')
$r.Value = "'" + '  -3.62%  '
$r.Style = "Normal"
$r = $ws.Range('E36')
$r.Value = "'" + '  -0.41%  '
$r.Style = "Normal"
$r = $ws.Range('D37')
$r.Value = "'" + '0.951'
$r.Style = "Normal"
$r = $ws.Range('E37')
$r.Value = "'" + '  +5.81%  '
$r.Style = "Normal"
$r = $ws.Range('D38')
$r.Value = "'" + '0.557'
$r.Style = "Normal"
$r = $ws.Range('E38')
$r.Value = "'" + '  -2.52%  '
$r.Style = "Normal"
$r = $ws.Range('E39')
$r.Value = "'" + '  -0.53%  '
$r.Style = "Normal"
$r = $ws.Range('D40')
$r.Value = "'" + '0.858'
$r.Style = "Normal"
$r = $ws.Range('E40')
$r.Value = "'" + '  -3.01%  '
$r.Style = "Normal"
$r = $ws.Range('E41')
$r.Value = "'" + '  +0.03%  '
$r.Style = "Normal"
$r = $ws.Range('D42')
$r.Value = "'" + '68.01'
$r.Style = "Normal"
$r = $ws.Range('E42')
$r.Value = "'" + '  +2.62%  '
$r.Style = "Normal"
$r = $ws.Range('E43')
$r.Value = "'" + '  +0.39%  '
$r.Style = "Normal"
$r = $ws.Range('D44')
$r.Value = "'" + '0.985'
$r.Style = "Normal"
$r = $ws.Range('E44')
$r.Value = "'" + '  -4.43%  '
$r.Style = "Normal"
$r = $ws.Range('E45')
$r.Value = "'" + '  -2.21%  '
$r.Style = "Normal"
$r = $ws.Range('D46')
$r.Value = "'" + '5.27'
$r.Style = "Normal"
$r = $ws.Range('E46')
$r.Value = "'" + '  -7.60%  '
$r.Style = "Normal"
$r = $ws.Range('D47')
$r.Value = "'" + '1.756.28'
$r.Style = "Normal"
$r = $ws.Range('E47')
$r.Value = "'" + '  -1.73%  '
$r.Style = "Normal"
$r = $ws.Range('E48')
$r.Value = "'" + '  +1.37%  '
$r.Style = "Normal"
$r = $ws.Range('D49')
$r.Value = "'" + '86.50'
$r.Style = "Normal"
$r = $ws.Range('E49')
$r.Value = "'" + '  -0.05%  '
$r.Style = "Normal"
$r = $ws.Range('D50')
$r.Value = "'" + '0.0₆0104'
$r.Style = "Normal"
$r = $ws.Range('E50')
$r.Value = "'" + '  -3.28%  '
$r.Style = "Normal"
$r = $ws.Range('E51')
$r.Value = "'" + '  +1.63%  '
$r.Style = "Normal"
